$wb = $excel.ActiveWorkbook

# Sheet 1: 土地
$ws = $wb.Worksheets.Item(1)
$ws.Range("B1").Value = 'name'
$ws.Range("C1").Value = 'area'
$ws.Range("D1").Value = 'share_portion'
$ws.Range("E1").Value = 'owner'
$ws.Range("F1").Value = 'register_date'
$ws.Range("G1").Value = 'register_reason'
$ws.Range("H1").Value = 'acquire_value'
$ws.Range("I1").Value = 'property_category'
$ws.Range("J1").Value = 'category'
$ws.Range("K1").Value = 'date'
$ws.Range("L1").Value = 'legislator_name'
$ws.Range("M1").Value = 'legislator_id'
$ws.Range("N1").Value = 'source_file'
$ws.Range("O1").Value = 'index'
$ws.Range("B2").Value = '臺北市松山區吳興段二小段05750000地號'
$ws.Range("D2").Value = '10000分之154'
$ws.Range("F2").Value = '75年10月30日'
$ws.Range("I2").Value = 'land'
$ws.Range("J2").Value = 'normal'
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = '2012-04-20'
$ws.Range("L2").Value = '費鴻泰'
$ws.Range("M2").Value = 1365
$ws.Range("N2").Value = 'tmpe52e1'
$ws.Range("O2").Value = 13
$ws.Range("B3").Value = '臺北市內湖區西湖段一小段00520012地號'
$ws.Range("D3").Value = '10000分之195'
$ws.Range("F3").Value = '79年09月15日'
$ws.Range("I3").Value = 'land'
$ws.Range("J3").Value = 'normal'
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = '2012-04-20'
$ws.Range("L3").Value = '費鴻泰'
$ws.Range("M3").Value = 1365
$ws.Range("N3").Value = 'tmpe52e1'
$ws.Range("O3").Value = 14
$ws.Range("B4").Value = '臺北市松山區延吉段三小段08400000地號'
$ws.Range("D4").Value = '85659分之7217'
$ws.Range("F4").Value = '87年12月10臼'
$ws.Range("I4").Value = 'land'
$ws.Range("J4").Value = 'normal'
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = '2012-04-20'
$ws.Range("L4").Value = '費鴻泰'
$ws.Range("M4").Value = 1365
$ws.Range("N4").Value = 'tmpe52e1'
$ws.Range("O4").Value = 15
$ws.Range("B5").Value = '臺北市內湖區西湖段四小段00410000地號'
$ws.Range("D5").Value = '20000分之263'
$ws.Range("F5").Value = '96年06月20日'
$ws.Range("G5").Value = 'ccdsK貝賣'
$ws.Range("I5").Value = 'land'
$ws.Range("J5").Value = 'normal'
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = '2012-04-20'
$ws.Range("L5").Value = '費鴻泰'
$ws.Range("M5").Value = 1365
$ws.Range("N5").Value = 'tmpe52e1'
$ws.Range("O5").Value = 16
$ws.Range("B6").Value = '臺北市內湖區西湖段四小段004]0000地號'
$ws.Range("D6").Value = '20000分之263'
$ws.Range("F6").Value = '96年06月20日'
$ws.Range("I6").Value = 'land'
$ws.Range("J6").Value = 'normal'
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = '2012-04-20'
$ws.Range("L6").Value = '費鴻泰'
$ws.Range("M6").Value = 1365
$ws.Range("N6").Value = 'tmpe52e1'
$ws.Range("O6").Value = 17

# Sheet 2: 建物
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = '臺北市松山區吳興段二小段01198000建號(4樓）'
$ws.Range("F2").Value = '75年10月30日'
$ws.Range("B3").Value = '臺北市內湖區西湖段小段04816000建號(4樓5樓）'
$ws.Range("F3").Value = '79年09月15H'
$ws.Range("B4").Value = '臺北市內湖區西湖段一小段04880000建號(地下室停車位)'
$ws.Range("F4").Value = '79年09月15H'
$ws.Range("B5").Value = '臺北市松山區延吉段三小段03504000建號(平台9.24平方公尺花台0.64平方公尺）'
$ws.Range("F5").Value = '87年12月10曰'
$ws.Range("B6").Value = '臺北市內湖區西湖段四小段01312000建號(6樓）'
$ws.Range("F6").Value = '96年06月20日'
$ws.Range("B7").Value = '臺北市內湖區西湖段四小段01312000建號(6樓）'
$ws.Range("F7").Value = '96年06月20□'
$ws.Range("B8").Value = '臺北市內湖區西湖段四小段01316000建號(地下層第三層'
$ws.Range("F8").Value = '96年06月20日'
$ws.Range("B9").Value = '臺北市内湖區西湖段四小段01316000建號(地下層第三層停車位2個）'
$ws.Range("F9").Value = '96年06月20日'

# Sheet 3: 汽車
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = '□產TEANA'
$ws.Range("E2").Value = '101年04月19曰'
$ws.Range("B3").Value = '日產BLURBIRD'
$ws.Range("E3").Value = '101年04月19曰'

# Sheet 6: 基金受益憑證
$ws = $wb.Worksheets.Item(6)
$ws.Range("D2").Value = '友邦證券投資信託股份有限公司'
$ws.Range("D3").Value = '友邦證券投資信託股份有限公司'
$ws.Range("D4").Value = '友邦證券投資信託股份有限公司'

# Sheet 7: 債務
$ws = $wb.Worksheets.Item(7)
$ws.Range("D2").Value = '國泰世華臺北市内湖區内湖路'
$ws.Range("F2").Value = '96年06月20日'
